# Automatische test-sync: 2025-06-24 20:44:50
# Appends a new logged e-mail (row 22) to the "Logs" sheet, extends the
# conditional formatting ranges to cover it, and bumps the matching
# "Retour / Terugbetaling" tally on the "Dashboard" sheet from 5 to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$a22 = "Ruilen van product"
$b22 = "mailmind.test@zohomail.eu"
$c22 = "Kan ik dit product ruilen voor een andere maat?`nSent using {0}"
$d22 = "Retour / Terugbetaling"
$e22 = "Beste klant,`nBedankt voor je bericht. Om je vraag over het ruilen van het product voor een andere maat te beantwoorden: Ja, dat is mogelijk. We hebben een ruilbeleid waarbij je het product kunt omruilen voor een andere maat, mits het product in de originele staat verkeert en binnen de gestelde termijn wordt geretourneerd. `nGraag ontvangen wij meer informatie over je bestelling, zoals het ordernummer en de gewenste maat, zodat we je verder kunnen helpen met het proces van ruilen. Aarzel niet om contact met ons op te nemen als je nog vragen hebt.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$f22 = "2025-06-24 20:44:49"
$g22 = "Ja"

$ws.Range("A22").Value = $a22
$ws.Range("B22").Value = $b22
$ws.Range("C22").Value = $c22
$ws.Range("D22").Value = $d22
$ws.Range("E22").Value = $e22
$ws.Range("F22").Value = $f22
$ws.Range("G22").Value = $g22

# Multi-line cell content (C22/E22) otherwise leaves the row pinned at a
# manually-sized height; auto-fit restores the sheet's normal "no explicit
# row height" state, matching every other row.
$ws.Rows.Item(22).AutoFit()

# The "Categorie"/"Beantwoord" columns carry conditional formatting over
# the data rows; stretch both ranges down to include the freshly added row.
$fcsD = $ws.Range("D2:D21").FormatConditions
$fcsD.Item(1).ModifyAppliesToRange($ws.Range("D2:D22"))

$fcsG = $ws.Range("G2:G21").FormatConditions
$fcsG.Item(1).ModifyAppliesToRange($ws.Range("G2:G22"))

# Dashboard tally for "Retour / Terugbetaling" goes from 5 to 6.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
